# CMPG323 EcoPower Logistics Data - "writing to excel file working"
#
# Narrative: the user switches to the Customers sheet, fills column F
# ("Test Result") with TRUE for every existing customer row, the fill
# handle overshoots by one blank row, and the selection is left on F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# Bring the Customers sheet to the front (flips tabSelected / activeTab).
$ws.Activate()

# Write a boolean TRUE down column F for rows 2-14 (the data rows) plus
# one extra row (15) picked up by the fill-handle overshoot.
$ws.Range("F2:F15").Value = $true

# Leave the selection where the user's last click landed.
[void]$ws.Range("F3").Select()
